# Applies the "Modificaciones finales segunda entrega" edit to the
# "diccionario de datos" workbook:
#  - removes the now-unused "Tipo de datos" column (C)
#  - makes the header row (Columna / Descripción) bold
#  - fills in the missing description for "ecomm_tipo_envio" (row 16)
#  - documents two new fields: sex / client_id (rows 35-36)
#  - leaves the last used cell (B16) selected, matching the author's session

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Document two additional columns that were added to the dataset.
$ws.Range("A35").Value = "sex"
$ws.Range("B35").Value = "Sexo del comprador"
$ws.Range("A36").Value = "client_id"
$ws.Range("B36").Value = "Id único del cliente"

# Add description for the existing "ecomm_tipo_envio" row that was blank.
$ws.Range("B16").Value = "Si es a domicilio o punto de retiro"

# The extra "Tipo de datos" header/column is no longer used; clear its content.
$ws.Range("C1").ClearContents()

# Highlight the header row in bold.
$ws.Range("A1:B1").Font.Bold = $true

# Restore the selection to match the author's last-saved cursor position.
$ws.Range("B16").Select()
